$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -58.2985

$ws.Range("B3").Value = -447.9618
$ws.Range("C3").Value = -60.6531

$ws.Range("C4").Value = -107.1259
$ws.Range("C5").Value = -176.8828
$ws.Range("C6").Value = -171.0812
$ws.Range("C7").Value = -116.6926
$ws.Range("C8").Value = -115.1218
$ws.Range("C9").Value = -129.6293
$ws.Range("C10").Value = -110.4751
$ws.Range("C11").Value = -78.3314
$ws.Range("C12").Value = -49.4337
$ws.Range("C13").Value = -50.9366
$ws.Range("C14").Value = -36.2095
$ws.Range("C15").Value = -49.2663
$ws.Range("C16").Value = -7.6823
$ws.Range("C17").Value = -1.2013
$ws.Range("C18").Value = -49.834
$ws.Range("C19").Value = -53.4183
$ws.Range("C20").Value = -87.08450000000001
$ws.Range("C21").Value = -45.1746
$ws.Range("C22").Value = -11.8291
$ws.Range("C23").Value = 15.9913
$ws.Range("C24").Value = 29.2182
